$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-10-30"

# Update the header label in I1 ("2022 (through 10-29)" -> "2022 (through 10-30)")
$ws.Range("I1").Value = "2022 (through 10-30)"

# Update the October and November values in column I (2022 total col)
$ws.Range("I10").Value = 144
$ws.Range("I11").Value = 121

# Update the Total row value
$ws.Range("I14").Value = 1397
